# ============================================================================
# PlayerPerformance_4536.xlsx update
#   - Insert a new "Player Info" sheet before "ODI Batting"
#   - Insert a new "ODI Batting Extra" sheet after "ODI Bowling"
#   - Rename MATCH_CARD_LINK -> MATCH_CODE on both existing sheets and
#     replace the full scorecard URL values with just the bare match code
#   - Clear a few stray empty INNING_NUMBER cells on "ODI Batting"
# ============================================================================

# Writes $value into a cell as TEXT, even when it looks numeric (match codes,
# ids, etc. are stored as text in this workbook). A leading apostrophe is the
# standard Excel idiom for "treat what follows as text".
function Set-TextCell($ws, $row, $col, $value) {
    if ($value -eq $null) {
        $ws.Cells.Item($row, $col).ClearContents()
        return
    }
    $ws.Cells.Item($row, $col).Value = "'" + $value
}

# Writes $value into a cell as a real NUMBER.
function Set-NumCell($ws, $row, $col, $value) {
    if ($value -eq $null) {
        $ws.Cells.Item($row, $col).ClearContents()
        return
    }
    $ws.Cells.Item($row, $col).Value = $value
}

# Bold header w/ thin border + center/top alignment, matching the look of
# the header row already used on "ODI Batting" / "ODI Bowling".
function Format-HeaderRange($range) {
    $range.Font.Bold = $true
    $range.Borders.LineStyle = 1
    $range.HorizontalAlignment = -4108
    $range.VerticalAlignment = -4160
}

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet layout: Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# ---------------------------------------------------------------------------
$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")

$playerInfoWs = $wb.Worksheets.Add($battingWs)
$playerInfoWs.Name = "Player Info"

# re-fetch ODI Bowling - earlier refs can go stale once the sheet collection
# is mutated
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")
$extraWs = $wb.Worksheets.Add($null, $bowlingWs)
$extraWs.Name = "ODI Batting Extra"

# re-fetch the two pre-existing sheets now that the collection has settled
$battingWs = $wb.Worksheets.Item("ODI Batting")
$bowlingWs = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------------
# 2. "Player Info" sheet content
# ---------------------------------------------------------------------------
$playerInfoHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $playerInfoHeaders.Length; $c++) {
    Set-TextCell $playerInfoWs 1 $c $playerInfoHeaders[$c - 1]
}
Format-HeaderRange $playerInfoWs.Range("A1:D1")

Set-TextCell $playerInfoWs 2 1 "4536"
Set-TextCell $playerInfoWs 2 2 "Mohammad Nawaz"
Set-TextCell $playerInfoWs 2 3 "Left Handed"
Set-TextCell $playerInfoWs 2 4 "Left Arm Orthodox"

# ---------------------------------------------------------------------------
# 3. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE (header + values), and
#    clear the stray empty INNING_NUMBER cells on rows 8/12/16/22/24
# ---------------------------------------------------------------------------
Set-TextCell $battingWs 1 4 "MATCH_CODE"

$battingCodes = @{
    2  = "3921"; 3  = "3925"; 4  = "3928"; 5  = "3930"; 6  = "3932";
    7  = "3939"; 8  = "3943"; 9  = "3944"; 10 = "3972"; 11 = "4114";
    12 = "4178"; 13 = "4200"; 14 = "4201"; 15 = "4204"; 16 = "4376";
    17 = "4460"; 18 = "4586"; 19 = "4590"; 20 = "4592"; 21 = "4634";
    22 = "4638"; 23 = "4641"; 24 = "4686"; 25 = "4688"; 26 = "4690";
}
foreach ($r in $battingCodes.Keys) {
    Set-TextCell $battingWs $r 4 $battingCodes[$r]
}

foreach ($r in @(8, 12, 16, 22, 24)) {
    $battingWs.Cells.Item($r, 2).ClearContents()
}

# ---------------------------------------------------------------------------
# 4. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE (header + values)
# ---------------------------------------------------------------------------
Set-TextCell $bowlingWs 1 2 "MATCH_CODE"

$bowlingCodes = @{
    2  = "3921"; 3  = "3925"; 4  = "3928"; 5  = "3930"; 6  = "3932";
    7  = "3939"; 8  = "3943"; 9  = "3944"; 10 = "3972"; 11 = "4114";
    12 = "4178"; 13 = "4200"; 14 = "4201"; 15 = "4204"; 16 = "4376";
    17 = "4460"; 18 = "4586"; 19 = "4590"; 20 = "4592"; 21 = "4634";
    22 = "4638"; 23 = "4641"; 24 = "4686"; 25 = "4688"; 26 = "4690";
}
foreach ($r in $bowlingCodes.Keys) {
    Set-TextCell $bowlingWs $r 2 $bowlingCodes[$r]
}

# ---------------------------------------------------------------------------
# 5. "ODI Batting Extra" sheet content
# ---------------------------------------------------------------------------
$extraHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $extraHeaders.Length; $c++) {
    Set-TextCell $extraWs 1 $c $extraHeaders[$c - 1]
}
Format-HeaderRange $extraWs.Range("A1:F1")

# MATCH_CODE, BATTING_POSITION(number|$null), NUM_4, NUM_6, PERCENT_RUNS_OF_TOTAL, MAN_OF_MATCH
$extraRows = @(
    @("3939", 8,    "0", "1", "6.69%",  "NO"),
    @("3943", 8,    $null, $null, $null, "NO"),
    @("3944", 8,    "0", "0", "1.30%",  "NO"),
    @("3972", 8,    "0", "0", "0.57%",  "NO"),
    @("4114", 10,   "3", "1", "8.98%",  "NO"),
    @("4178", $null,$null, $null, $null, "NO"),
    @("4200", 8,    "0", "1", "3.88%",  "NO"),
    @("4201", 8,    "1", "0", "6.33%",  "NO"),
    @("4204", 8,    "0", "0", "3.96%",  "NO"),
    @("4376", 8,    $null, $null, $null, "NO"),
    @("4460", 7,    "0", "0", "1.25%",  "NO"),
    @("4586", $null,$null, $null, $null, "NO"),
    @("4590", $null,$null, $null, $null, "NO"),
    @("4592", 8,    "0", "0", "2.97%",  "NO"),
    @("4634", $null,$null, $null, $null, "NO"),
    @("4638", 8,    $null, $null, $null, "YES"),
    @("4641", 7,    "1", "2", "13.11%", "NO"),
    @("4686", $null,$null, $null, $null, "NO"),
    @("4688", 7,    "0", "0", "1.65%",  "NO"),
    @("4690", $null,$null, $null, $null, "NO")
)

$r = 2
foreach ($row in $extraRows) {
    Set-TextCell $extraWs $r 1 $row[0]
    Set-NumCell  $extraWs $r 2 $row[1]
    Set-TextCell $extraWs $r 3 $row[2]
    Set-TextCell $extraWs $r 4 $row[3]
    Set-TextCell $extraWs $r 5 $row[4]
    Set-TextCell $extraWs $r 6 $row[5]
    $r++
}

Write-Host "done"
